$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.526.06"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "'1.833.43"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.73%  "

$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("D7").Value = "'0.4595"

$ws.Range("D8").Value = "'0.3828"
$ws.Range("E8").Value = "  -1.70%  "

$ws.Range("D9").Value = "'46.33"
$ws.Range("E9").Value = "  +1.43%  "

$ws.Range("D10").Value = "'0.07849"
$ws.Range("E10").Value = "  -1.08%  "

$ws.Range("D11").Value = "'0.9612"
$ws.Range("E11").Value = "  -4.22%  "

$ws.Range("D12").Value = "'21.08"
$ws.Range("E12").Value = "  -2.41%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.848"
$ws.Range("E13").Value = "  -1.87%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.816.12"
$ws.Range("E14").Value = "  -1.56%  "

$ws.Range("D15").Value = "'7.077"
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  -0.77%  "

$ws.Range("D17").Value = "'89.69"
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("D18").Value = "'0.06570"
$ws.Range("E18").Value = "  -2.13%  "

$ws.Range("D19").Value = "'0.00001021"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").Value = "'17.13"
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "  -0.66%  "

$ws.Range("D22").Value = "'27.486.03"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").Value = "'5.304"
$ws.Range("E23").Value = "  -2.11%  "

$ws.Range("D24").Value = "'10.79"
$ws.Range("E24").Value = "  -1.21%  "

$ws.Range("D25").Value = "'2.271"
$ws.Range("E25").Value = "  -1.63%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'158.72"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "'2.041.42"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("D28").Value = "'19.38"
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("D29").Value = "'2.043"
$ws.Range("E29").Value = "  -4.23%  "

$ws.Range("D30").Value = "'5.287"
$ws.Range("E30").Value = "  -3.00%  "

$ws.Range("D31").Value = "'118.07"
$ws.Range("E31").Value = "  -2.76%  "

$ws.Range("D32").Value = "'0.09400"

$ws.Range("D33").Value = "'0.9313"
$ws.Range("E33").Value = "  -4.49%  "

$ws.Range("D34").Value = "'3.573"
$ws.Range("E34").Value = "  -1.45%  "

$ws.Range("D35").Value = "'5.212"
$ws.Range("E35").Value = "  -2.19%  "

$ws.Range("D36").Value = "'1.318"
$ws.Range("E36").Value = "  -1.71%  "

$ws.Range("D37").Value = "'0.05946"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").Value = "'0.02185"
$ws.Range("E38").Value = "  -2.11%  "

$ws.Range("D39").Value = "'8.133"
$ws.Range("E39").Value = "  -2.73%  "

$ws.Range("D40").Value = "'1.004"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("D41").Value = "'1.144"
$ws.Range("E41").Value = "  -3.82%  "

$ws.Range("D42").Value = "'0.5737"
$ws.Range("E42").Value = "  -3.28%  "

$ws.Range("D43").Value = "'0.1824"
$ws.Range("E43").Value = "  -2.18%  "

$ws.Range("D44").Value = "'9.959"
$ws.Range("E44").Value = "  -4.60%  "

$ws.Range("D45").Value = "'1.266"
$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5384"
$ws.Range("E46").Value = "  -3.68%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'11.77"
$ws.Range("E47").Value = "  -3.45%  "

$ws.Range("D48").Value = "'1.904"
$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("D49").Value = "'0.06824"
$ws.Range("E49").Value = "  +1.63%  "

$ws.Range("D50").Value = "'111.05"
$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  -32.72%  "
